$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.975.64"
$ws.Range("E2").Value = "  +3.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.241.30"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.86"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "80.17"
$ws.Range("E6").Value = "  +8.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  +2.34%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.601"
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.23"
$ws.Range("E10").Value = "  +8.01%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.07"
$ws.Range("E12").Value = "  +3.80%  "
$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.572.88"
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.72"
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.266.43"
$ws.Range("E16").Value = "  +3.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.791"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.886.49"
$ws.Range("E18").Value = "  +3.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000104"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.30"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  +2.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.36"
$ws.Range("E22").Value = "  +7.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.33"
$ws.Range("E23").Value = "  +2.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.40"
$ws.Range("E24").Value = "  -2.97%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.86"
$ws.Range("E26").Value = "  +2.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "40.50"
$ws.Range("E27").Value = "  +9.37%  "
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("E30").Value = "  +3.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.39"
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("E32").Value = "  +10.55%  "
$ws.Range("E33").Value = "  +2.79%  "
$ws.Range("E34").Value = "  +2.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.115"
$ws.Range("E35").Value = "  +7.22%  "
$ws.Range("E36").Value = "  +2.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0370"
$ws.Range("E37").Value = "  +11.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.53"
$ws.Range("E38").Value = "  +4.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.14"
$ws.Range("E39").Value = "  +8.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.99"
$ws.Range("E40").Value = "  +24.29%  "
$ws.Range("E41").Value = "  +3.20%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.43"
$ws.Range("E42").Value = "  +7.46%  "
$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.52"
$ws.Range("E43").Value = "  +5.95%  "
$ws.Range("E44").Value = "  +2.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.43"
$ws.Range("E45").Value = "  +1.76%  "
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0987"
$ws.Range("E47").Value = "  +1.34%  "
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.55"
$ws.Range("E49").Value = "  +28.37%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.447"
$ws.Range("E50").Value = "  -4.64%  "
$ws.Range("E51").Value = "  +2.89%  "
